$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the Kannada, Hindi and Tamil language rows entirely.
$ws.Range("A5:E7").EntireRow.Delete() | Out-Null

# The Arabic row's native_name ("Arabic") was wrong - it should be the
# Arabic word for Arabic. Clear the row back to the sheet's default style
# and re-enter the (now corrected) values.
$ws.Range("A4:D4").Clear() | Out-Null
$ws.Range("A4").Value = "ara"
$ws.Range("B4").Value = "Arabic"
$ws.Range("C4").Value = "الهندو أوروبية"
$ws.Range("D4").Value = "عربي"

# Leave the selection where the user ended up after making the edit.
$ws.Range("F11").Select() | Out-Null
